$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Url   = "https://www.genomeweb.com/cancer/guardant-health-blood-test-receives-japanese-approval-cdx-eli-lilly-breast-cancer-drug"
        Kw    = "CDx"
        Title = "Guardant Health Blood Test Receives Japanese Approval as CDx for Eli Lilly Breast Cancer Drug"
    },
    @{
        Url   = "https://www.360dx.com/cancer/guardant-health-blood-test-receives-japanese-approval-cdx-eli-lilly-breast-cancer-drug"
        Kw    = "CDx"
        Title = "Guardant Health Blood Test Receives Japanese Approval as CDx for Eli Lilly Breast Cancer Drug"
    }
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $a = $ws.Cells.Item($r, 1)
    $a.Value2 = $row.Url
    $ws.Hyperlinks.Add($a, $row.Url) | Out-Null
    $a.Style = "Hyperlink"

    $ws.Cells.Item($r, 2).Value2 = $row.Kw
    $ws.Cells.Item($r, 3).Value2 = $row.Title
}
